$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H101").Value = 334.5
$ws.Range("I101").Value = 334.5
$ws.Range("K101").Value = 1003.5
$ws.Range("M101").Value = 618.5
$ws.Range("H112").Value = 2150.5
$ws.Range("J112").Value = 2210.182
$ws.Range("L112").Value = 6630.545999999999
$ws.Range("N112").Value = -8846.545999999998
$ws.Range("H123").Value = 73240
$ws.Range("J123").Value = 73240
$ws.Range("L123").Value = 73240
$ws.Range("N123").Value = -83040
$ws.Range("H132").Value = 12670.886
$ws.Range("I132").Value = 12928.677
$ws.Range("K132").Value = 38786.031
$ws.Range("M132").Value = -36256.031
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 78950
$ws.Range("J140").Value = 78950
$ws.Range("L140").Value = 78950
$ws.Range("N140").Value = -89310
$ws.Range("H141").Value = 12339.96
$ws.Range("I141").Value = 4374.75
$ws.Range("J141").Value = 13857.143
$ws.Range("K141").Value = 13124.25
$ws.Range("L141").Value = 41571.429
$ws.Range("M141").Value = -7944.25
$ws.Range("N141").Value = -51931.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1104.0834
$ws.Range("I2").Value = 995.2857
$ws.Range("K2").Value = 995.2857
$ws.Range("M2").Value = -882.2857
$ws.Range("H61").Value = 3848236.2
$ws.Range("I61").Value = 2242.9092
$ws.Range("K61").Value = 2242.9092
$ws.Range("M61").Value = -2030.9092
$ws.Range("H74").Value = 1042410.75
$ws.Range("I74").Value = 1178952.9
$ws.Range("J74").Value = 31999.2
$ws.Range("K74").Value = 1178952.9
$ws.Range("L74").Value = 31999.2
$ws.Range("M74").Value = -1178078.9
$ws.Range("N74").Value = -33747.2
$ws.Range("H77").Value = 1042410.75
$ws.Range("I77").Value = 1178952.9
$ws.Range("J77").Value = 31999.2
$ws.Range("K77").Value = 5894764.5
$ws.Range("L77").Value = 159996
$ws.Range("M77").Value = -5890396.5
$ws.Range("N77").Value = -168732
$ws.Range("H102").Value = 3848.9
$ws.Range("I102").Value = 3882.1667
$ws.Range("K102").Value = 3882.1667
$ws.Range("M102").Value = -2260.1667
$ws.Range("H110").Value = 1945.5555
$ws.Range("I110").Value = 1833.3334
$ws.Range("K110").Value = 1833.3334
$ws.Range("M110").Value = 211.6666
$ws.Range("H116").Value = 1104.0834
$ws.Range("I116").Value = 995.2857
$ws.Range("K116").Value = 995.2857
$ws.Range("M116").Value = 1298.7143
$ws.Range("H132").Value = 4281.2666
$ws.Range("I132").Value = 2097.0527
$ws.Range("J132").Value = 8054
$ws.Range("K132").Value = 6291.158100000001
$ws.Range("L132").Value = 24162
$ws.Range("M132").Value = -3761.158100000001
$ws.Range("N132").Value = -29222
$ws.Range("H136").Value = 3848236.2
$ws.Range("I136").Value = 2242.9092
$ws.Range("K136").Value = 6728.7276
$ws.Range("M136").Value = -4178.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1104.0834
$ws.Range("I3").Value = 995.2857
$ws.Range("K3").Value = 995.2857
$ws.Range("M3").Value = -881.2857
$ws.Range("H22").Value = 297.75
$ws.Range("I22").Value = 297.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 297.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -124.75
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1911.52
$ws.Range("I132").Value = 1794.4524
$ws.Range("J132").Value = 2526.125
$ws.Range("K132").Value = 5383.357199999999
$ws.Range("L132").Value = 7578.375
$ws.Range("M132").Value = -2853.357199999999
$ws.Range("N132").Value = -12638.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 868.1818
$ws.Range("J92").Value = 1850
$ws.Range("L92").Value = 5550
$ws.Range("N92").Value = -8046
$ws.Range("H129").Value = 2681
$ws.Range("I129").Value = 2500
$ws.Range("J129").Value = 2696.0833
$ws.Range("K129").Value = 7500
$ws.Range("L129").Value = 8088.249899999999
$ws.Range("M129").Value = -2500
$ws.Range("N129").Value = -18088.2499
$ws.Range("H131").Value = 6192.5806
$ws.Range("I131").Value = 1248.25
$ws.Range("K131").Value = 3744.75
$ws.Range("M131").Value = 1295.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.666664
$ws.Range("I2").Value = 59.214287
$ws.Range("K2").Value = 59.214287
$ws.Range("M2").Value = 53.785713
$ws.Range("H47").Value = 16500
$ws.Range("I47").Value = 15000
$ws.Range("K47").Value = 15000
$ws.Range("M47").Value = -14432
$ws.Range("H70").Value = 10523.088
$ws.Range("I70").Value = 11680.889
$ws.Range("J70").Value = 6057.2856
$ws.Range("K70").Value = 11680.889
$ws.Range("L70").Value = 6057.2856
$ws.Range("M70").Value = -11410.889
$ws.Range("N70").Value = -6597.2856
$ws.Range("H73").Value = 10523.088
$ws.Range("I73").Value = 11680.889
$ws.Range("J73").Value = 6057.2856
$ws.Range("K73").Value = 11680.889
$ws.Range("L73").Value = 6057.2856
$ws.Range("M73").Value = -10744.889
$ws.Range("N73").Value = -7929.2856
$ws.Range("H107").Value = 3816.3333
$ws.Range("I107").Value = 3224.75
$ws.Range("J107").Value = 4999.5
$ws.Range("K107").Value = 3224.75
$ws.Range("L107").Value = 4999.5
$ws.Range("M107").Value = -1304.75
$ws.Range("N107").Value = -8839.5
$ws.Range("H132").Value = 14320.593
$ws.Range("I132").Value = 8031.7617
$ws.Range("K132").Value = 24095.2851
$ws.Range("M132").Value = -21565.2851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2911.524
$ws.Range("I22").Value = 2474.6
$ws.Range("K22").Value = 2474.6
$ws.Range("M22").Value = -2179.6
$ws.Range("H27").Value = 2911.524
$ws.Range("I27").Value = 2474.6
$ws.Range("K27").Value = 2474.6
$ws.Range("M27").Value = -2367.6
$ws.Range("H40").Value = 5006.778
$ws.Range("I40").Value = 3735.1667
$ws.Range("K40").Value = 3735.1667
$ws.Range("M40").Value = -3599.1667
$ws.Range("H55").Value = 1393.4103
$ws.Range("I55").Value = 1151.6666
$ws.Range("J55").Value = 1600.619
$ws.Range("K55").Value = 1151.6666
$ws.Range("L55").Value = 1600.619
$ws.Range("M55").Value = -978.6666
$ws.Range("N55").Value = -1946.619

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 37156.312
$ws.Range("I122").Value = 1903.1852
$ws.Range("J122").Value = 227523.2
$ws.Range("K122").Value = 5709.5556
$ws.Range("L122").Value = 682569.6000000001
$ws.Range("M122").Value = -3259.5556
$ws.Range("N122").Value = -687469.6000000001
$ws.Range("H132").Value = 4067753.8
$ws.Range("I132").Value = 4275869.5
$ws.Range("J132").Value = 9502
$ws.Range("K132").Value = 12827608.5
$ws.Range("L132").Value = 28506
$ws.Range("M132").Value = -12825078.5
$ws.Range("N132").Value = -33566
$ws.Range("H136").Value = 5096485
$ws.Range("I136").Value = 2719701.2
$ws.Range("K136").Value = 8159103.600000001
$ws.Range("M136").Value = -8156553.600000001
$ws.Range("H139").Value = 82000
$ws.Range("J139").Value = 82000
$ws.Range("L139").Value = 82000
$ws.Range("N139").Value = -92280
$ws.Range("H140").Value = 47595.4
$ws.Range("J140").Value = 47595.4
$ws.Range("L140").Value = 47595.4
$ws.Range("N140").Value = -57955.4
$ws.Range("H141").Value = 69933
$ws.Range("J141").Value = 69933
$ws.Range("L141").Value = 69933
